# edit.ps1 — apply the "add 2022-Q3 data" commit to the workbook.
#
# Summary of the change:
#   * A new worksheet "2022-Q3" is inserted right after "总计" (so it
#     becomes the second tab); every quarter sheet after it shifts right
#     by one position. The new sheet holds the same 7-column fund-holding
#     table layout used by the other quarter sheets.
#   * The "总计" (summary) sheet gets a new top data row for 2022-Q3 and
#     all the other quarters shift down by one row.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "总计" sheet — shift existing quarters down a row and insert the
#    new 2022-Q3 figures at the top of the table (row 2).
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Make room for a new row 9 (2020-Q4 moves here); clone formatting from A8
$ws1.Range("A8").Copy()
$ws1.Range("A9").PasteSpecial(-4122)

# Row 9 <= old row 8 (2020-Q4)
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "2020-Q4"
$ws1.Range("C9").Value = 10
$ws1.Range("D9").Value = 11.16

# Row 8 <= old row 7 (2021-Q1)
$ws1.Range("B8").Value = "2021-Q1"
$ws1.Range("C8").Value = 17
$ws1.Range("D8").Value = 16.4

# Row 7 <= old row 6 (2021-Q2)
$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 19
$ws1.Range("D7").Value = 14.88

# Row 6 <= old row 5 (2021-Q3)
$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 39
$ws1.Range("D6").Value = 27.37

# Row 5 <= old row 4 (2021-Q4)
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 62
$ws1.Range("D5").Value = 38.47

# Row 4 <= old row 3 (2022-Q1)
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 55
$ws1.Range("D4").Value = 32.02

# Row 3 <= old row 2 (2022-Q2)
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 19
$ws1.Range("D3").Value = 21.69

# Row 2 <= brand-new 2022-Q3 figures
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 16
$ws1.Range("D2").Value = 17.62

# -----------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计".
# -----------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q3"

# Header row (B1:H1) — reuse the bold/centered/bordered header style
# used by every other quarter sheet (copy format from "总计"'s header).
$ws1.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows 2..17 — column A (index) and H (rank) are numeric; columns
# B..G are stored as text (same as every other quarter sheet), including
# the numeric-looking fund codes / sizes / percentages.
$dataRange = $newSheet.Range("B2:G17")
$dataRange.NumberFormat = "@"

# Row 2: 163406
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "163406"
$newSheet.Range("C2").Value = "兴全合润混合"
$newSheet.Range("D2").Value = "252.62"
$newSheet.Range("E2").Value = "92.30"
$newSheet.Range("F2").Value = "3.23"
$newSheet.Range("G2").Value = "8.1596"
$newSheet.Range("H2").Value = 9

# Row 3: 163417
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "163417"
$newSheet.Range("C3").Value = "兴全合宜灵活配置混合（LOF）A"
$newSheet.Range("D3").Value = "163.32"
$newSheet.Range("E3").Value = "92.65"
$newSheet.Range("F3").Value = "3.29"
$newSheet.Range("G3").Value = "5.3732"
$newSheet.Range("H3").Value = 7

# Row 4: 519692
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "519692"
$newSheet.Range("C4").Value = "交银成长混合A"
$newSheet.Range("D4").Value = "23.45"
$newSheet.Range("E4").Value = "76.71"
$newSheet.Range("F4").Value = "7.01"
$newSheet.Range("G4").Value = "1.6438"
$newSheet.Range("H4").Value = 5

# Row 5: 519694
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "519694"
$newSheet.Range("C5").Value = "交银蓝筹混合"
$newSheet.Range("D5").Value = "16.10"
$newSheet.Range("E5").Value = "78.09"
$newSheet.Range("F5").Value = "6.98"
$newSheet.Range("G5").Value = "1.1238"
$newSheet.Range("H5").Value = 5

# Row 6: 005123
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "005123"
$newSheet.Range("C6").Value = "南方优享分红灵活配置混合A"
$newSheet.Range("D6").Value = "7.37"
$newSheet.Range("E6").Value = "92.25"
$newSheet.Range("F6").Value = "5.01"
$newSheet.Range("G6").Value = "0.3692"
$newSheet.Range("H6").Value = 8

# Row 7: 005491
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "005491"
$newSheet.Range("C7").Value = "兴全合宜灵活配置混合（LOF）C"
$newSheet.Range("D7").Value = "10.59"
$newSheet.Range("E7").Value = "92.65"
$newSheet.Range("F7").Value = "3.29"
$newSheet.Range("G7").Value = "0.3484"
$newSheet.Range("H7").Value = 7

# Row 8: 001208
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "001208"
$newSheet.Range("C8").Value = "诺安低碳经济股票A"
$newSheet.Range("D8").Value = "15.05"
$newSheet.Range("E8").Value = "86.13"
$newSheet.Range("F8").Value = "1.58"
$newSheet.Range("G8").Value = "0.2378"
$newSheet.Range("H8").Value = 9

# Row 9: 001551
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "001551"
$newSheet.Range("C9").Value = "天弘中证医药100指数型发起式 C"
$newSheet.Range("D9").Value = "8.58"
$newSheet.Range("E9").Value = "95.24"
$newSheet.Range("F9").Value = "1.49"
$newSheet.Range("G9").Value = "0.1278"
$newSheet.Range("H9").Value = 3

# Row 10: 001550
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "001550"
$newSheet.Range("C10").Value = "天弘中证医药100指数型发起式 A"
$newSheet.Range("D10").Value = "5.31"
$newSheet.Range("E10").Value = "95.24"
$newSheet.Range("F10").Value = "1.49"
$newSheet.Range("G10").Value = "0.0791"
$newSheet.Range("H10").Value = 3

# Row 11: 010349
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "010349"
$newSheet.Range("C11").Value = "诺安低碳经济股票C"
$newSheet.Range("D11").Value = "3.83"
$newSheet.Range("E11").Value = "86.13"
$newSheet.Range("F11").Value = "1.58"
$newSheet.Range("G11").Value = "0.0605"
$newSheet.Range("H11").Value = 9

# Row 12: 320020
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "320020"
$newSheet.Range("C12").Value = "诺安策略精选股票"
$newSheet.Range("D12").Value = "2.07"
$newSheet.Range("E12").Value = "84.26"
$newSheet.Range("F12").Value = "1.61"
$newSheet.Range("G12").Value = "0.0333"
$newSheet.Range("H12").Value = 8

# Row 13: 013441
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "013441"
$newSheet.Range("C13").Value = "西藏东财创新医疗六个月定开混合"
$newSheet.Range("D13").Value = "0.49"
$newSheet.Range("E13").Value = "82.53"
$newSheet.Range("F13").Value = "5.33"
$newSheet.Range("G13").Value = "0.0261"
$newSheet.Range("H13").Value = 5

# Row 14: 006587
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "006587"
$newSheet.Range("C14").Value = "南方优享分红灵活配置混合C"
$newSheet.Range("D14").Value = "0.46"
$newSheet.Range("E14").Value = "92.25"
$newSheet.Range("F14").Value = "5.01"
$newSheet.Range("G14").Value = "0.0230"
$newSheet.Range("H14").Value = 8

# Row 15: 960016
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "960016"
$newSheet.Range("C15").Value = "交银成长混合H"
$newSheet.Range("D15").Value = "0.16"
$newSheet.Range("E15").Value = "76.71"
$newSheet.Range("F15").Value = "7.01"
$newSheet.Range("G15").Value = "0.0112"
$newSheet.Range("H15").Value = 5

# Row 16: 008444
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "008444"
$newSheet.Range("C16").Value = "九泰动态策略灵活配置混合C"
$newSheet.Range("D16").Value = "0.06"
$newSheet.Range("E16").Value = "46.24"
$newSheet.Range("F16").Value = "3.27"
$newSheet.Range("G16").Value = "0.0020"
$newSheet.Range("H16").Value = 6

# Row 17: 008443
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "008443"
$newSheet.Range("C17").Value = "九泰动态策略灵活配置混合A"
$newSheet.Range("D17").Value = "0.04"
$newSheet.Range("E17").Value = "46.24"
$newSheet.Range("F17").Value = "3.27"
$newSheet.Range("G17").Value = "0.0013"
$newSheet.Range("H17").Value = 6

# Drop the temporary text NumberFormat so B:G carry no explicit style
# (matches the plain/unstyled data cells used on the other sheets);
# column A keeps the bold/border style copied onto it above.
$dataRange.Style = "Normal"

# Re-apply the bold/border "总计"-style formatting to column A (index)
# cells, same treatment as every other quarter sheet. PasteSpecial with
# xlPasteFormats (-4122) only touches formatting, so the index values set
# above are left untouched.
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A17").PasteSpecial(-4122)

Write-Output "2022-Q3 sheet inserted; 总计 updated"
